# "updates from pascale and feather team"
#
# Updates the `attribute` sheet of the microhabitat metadata workbook:
#  - revises several attribute_definition cells in column B with expanded /
#    clarified text (species list, substrate size ranges, woody cover
#    definitions, and a more detailed geomorphic-unit description that
#    folds in the reviewer's threaded-comment feedback)
#  - removes the two threaded comments on B8 and B29 now that their
#    feedback has been incorporated into the text
#  - grows the row height for B8/B29 to fit the now-longer wrapped text
#  - adjusts the sheet's zoom level and active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Resolve reviewer threaded comments on B8 and B29, then remove them
#    (their content has been folded into the revised definitions below).
# ---------------------------------------------------------------------
if ($ws.Range("B8").Comment -ne $null) {
    $ws.Range("B8").Comment.Delete()
}
if ($ws.Range("B29").Comment -ne $null) {
    $ws.Range("B29").Comment.Delete()
}

# ---------------------------------------------------------------------
# 2. Update attribute_definition text in column B.
# ---------------------------------------------------------------------

$species = @'
Fish species. Species = c("chinook salmon", "sacramento pikeminnow", "speckled dace", 
"steelhead trout (wild)", "steelhead trout, (clipped)", "tule perch")
'@
$ws.Range("B8").Value2 = $species

$ws.Range("B15").Value2 = "Percentage of  fine substrate within microhabitat plot (less than 0.05mm)"
$ws.Range("B16").Value2 = "Percentage of  sand substrate within microhabitat plot (0.05-2mm)"
$ws.Range("B17").Value2 = "Percentage of small gravel substrate within microhabitat plot (2-50mm)"
$ws.Range("B18").Value2 = "Percentage of large gravel substrate within microhabitat plot (50-150mm)"
$ws.Range("B19").Value2 = "Percentage of cobble substrate within microhabitat plot (150-300mm)"
$ws.Range("B20").Value2 = "Percentage of boulder substrate within microhabitat plot (greater than 300mm)"

$ws.Range("B22").Value2 = "Percentage of small woody cover within microhabitat plot. Small wood is generally considered to be less than 10 cm diameter at breast height."
$ws.Range("B23").Value2 = "Percentage of large woody cover within microhabitat plot. Large wood is generally considered to be more than 10 cm at breat height and consists of tree trunks or a large branch."

$geomorphic = @'
Geomorphic features of stream designated through visual observation. Features include = c("glide", "glide margin", "pool", "riffle", "riffle margin", "backwater"). A glide is characterized by its smooth, uniform flow and relatively shallow depth. A riffle is a shallow section where the flow of water is faster and more turbulent. A pool is characterized by deeper depth, slower flow, and relatively calm water. A backwater area is characterized by stagnant or slow-moving water often found in area where the main flow is obstructed. Riffle or glide margins are transitions between these geomorphic units.
'@
$ws.Range("B29").Value2 = $geomorphic

# ---------------------------------------------------------------------
# 3. Grow the row heights for the two rows whose text got much longer so
#    the wrapped text stays visible.
# ---------------------------------------------------------------------
$ws.Rows.Item(8).RowHeight = 40
$ws.Rows.Item(29).RowHeight = 60

# ---------------------------------------------------------------------
# 4. Update the view: zoom level and active selection.
# ---------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.Zoom = 170
[void]$ws.Range("B16").Select()
